# Rename the single worksheet to reflect the new "updated on" date.
# Excel automatically keeps the RIS_systemer defined name (which refers to
# 'Opdateret d. 02-12-2025'!$A$1:$G$14) in sync with the sheet's new name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Opdateret d. 05-12-2025"
